$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Update the translation keys from .sad. to .mtp.
$ws.Range("A2").Value = "inspection_needed_export.mtp.heading"
$ws.Range("A3").Value = "inspection_needed_import.mtp.heading"

# Select A3 as the active cell
$ws.Range("A3").Select()
